$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "reviews_count" column (E) is no longer populated by the scraper,
# so drop it entirely and shift the remaining columns (reviews_average,
# latitude, longitude, is_permanently_closed, gmaps_link,
# latest_review_date) one position to the left.
$ws.Columns.Item(5).Delete()
